$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style used to restore cells to their original (unstyled) appearance
# after forcing a text number format, so no stray style index is left on them.
$plainStyle = $ws.Range("C2").Style

$ws.Range("D2").Value = '65.664.72'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '3.444.19'
$ws.Range("E3").Value = '  -3.24%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.90'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  -2.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.13'
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = '  -6.80%  '
$ws.Range("D7").Value = '3.442.14'
$ws.Range("E7").Value = '  -3.21%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.35'
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = '  -6.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.122'
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  -8.79%  '
$ws.Range("E12").Value = '  -7.35%  '
$ws.Range("D13").Value = '4.026.21'
$ws.Range("E13").Value = '  -3.28%  '
$ws.Range("E14").Value = '  -9.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.40'
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  -9.62%  '
$ws.Range("D16").Value = '3.443.53'
$ws.Range("E16").Value = '  -3.42%  '
$ws.Range("D17").Value = '65.622.39'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("E18").Value = '  -2.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.90'
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = '  -10.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.87'
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = '  -6.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.72'
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = '  -7.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '393.60'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  -6.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.554'
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = '  -8.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.42'
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = '  -5.77%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = '3.586.70'
$ws.Range("E26").Value = '  -3.06%  '
$ws.Range("E27").Value = '  -9.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.23'
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = '  -9.06%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.25'
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = '  -9.39%  '
$ws.Range("E31").Value = '  -10.31%  '
$ws.Range("D32").Value = '3.451.47'
$ws.Range("E32").Value = '  -2.93%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  -5.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.01'
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = '  -6.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '172.42'
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.96'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  -9.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.19'
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = '  -10.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = '  -8.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.82'
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = '  -9.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0767'
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  -7.36%  '
$ws.Range("E42").Value = '  -4.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.71'
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = '  -4.39%  '
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("E45").Value = '  -13.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.62'
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  -11.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.93'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.60'
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = '  -7.35%  '
$ws.Range("E50").Value = '  -14.99%  '
$ws.Range("D51").Value = '2.215.46'
$ws.Range("E51").Value = '  -7.33%  '
